# "implement stack adt, upload notes"
# Adds four new sorting-algorithm rows (Merge/Quick/Radix/Counting Sort) to
# the OVERVIEW sheet, and restores the cursor/selection state on each sheet.

$wb = $excel.ActiveWorkbook

# --- OVERVIEW: append Merge Sort / Quick Sort / Radix Sort / Counting Sort ---
$overview = $wb.Worksheets.Item("OVERVIEW")

$newAlgorithms = @("Merge Sort", "Quick Sort", "Radix Sort", "Counting Sort")
$row = 5
foreach ($name in $newAlgorithms) {
    $overview.Range("A$row").Value = $name
    $overview.Range("B$row").Value = "MEDIUM"
    $overview.Range("C$row").Value = 912
    $row++
}

# --- restore per-sheet selections ---
$selectionSort = $wb.Worksheets.Item("SelectionSort")
$selectionSort.Activate()
$selectionSort.Range("A40").Select() | Out-Null

$insertionSort = $wb.Worksheets.Item("InsertionSort")
$insertionSort.Activate()
$insertionSort.Range("A43").Select() | Out-Null

$overview.Activate()
$overview.Range("D5").Select() | Out-Null

Write-Host "OVERVIEW rows added; selections restored"
